# Daily refresh of the cryptos list (GitHub Actions scheduled run).
# Updates the Price (D) and Volume(1h) (E) columns for each coin row, and
# swaps the Kaspa / VeChain rows (42/43), which changed rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "70.460.18"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.643.25"
$ws.Range("E3").Value = "  +5.04%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "576.82"
$ws.Range("E5").Value = "  -1.78%  "

# Row 6 - Solana
$ws.Range("D6").Value = "175.63"
$ws.Range("E6").Value = "  -2.15%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.634.52"
$ws.Range("E7").Value = "  +5.03%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  +1.53%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.19%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.196"
$ws.Range("E10").Value = "  -4.52%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "6.78"
$ws.Range("E11").Value = "  +23.54%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.603"
$ws.Range("E12").Value = "  +1.62%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "48.48"
$ws.Range("E13").Value = "  -1.83%  "

# Row 14 - ShibaInu
$ws.Range("D14").Value = "0.0000285"
$ws.Range("E14").Value = "  -0.64%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.233.19"
$ws.Range("E15").Value = "  +5.32%  "

# Row 16 - BitcoinCash
$ws.Range("D16").Value = "669.39"
$ws.Range("E16").Value = "  -3.30%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "8.84"
$ws.Range("E17").Value = "  +0.57%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.654.96"
$ws.Range("E18").Value = "  +5.45%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "70.597.81"
$ws.Range("E19").Value = "  -0.10%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.48%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "17.73"
$ws.Range("E21").Value = "  -0.74%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "11.37"
$ws.Range("E22").Value = "  -1.44%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.928"
$ws.Range("E23").Value = "  +1.74%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "17.07"

# Row 25 - Litecoin
$ws.Range("D25").Value = "100.23"
$ws.Range("E25").Value = "  -1.40%  "

# Row 26 - PancakeSwap
$ws.Range("D26").Value = "3.89"
$ws.Range("E26").Value = "  -1.81%  "

# Row 27 - ImmutableX
$ws.Range("D27").Value = "2.77"
$ws.Range("E27").Value = "  +2.43%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  -0.02%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "9.94"
$ws.Range("E29").Value = "  +2.25%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "34.82"
$ws.Range("E30").Value = "  +2.39%  "

# Row 31 - Stacks
$ws.Range("D31").Value = "3.34"
$ws.Range("E31").Value = "  -0.54%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "8.95"
$ws.Range("E32").Value = "  +0.66%  "

# Row 33 - Mantle
$ws.Range("E33").Value = "  -5.33%  "

# Row 34 - NEARProtocol
$ws.Range("D34").Value = "7.28"
$ws.Range("E34").Value = "  +0.58%  "

# Row 35 - dogwifhat
$ws.Range("D35").Value = "3.98"
$ws.Range("E35").Value = "  +0.76%  "

# Row 36 - Bittensor
$ws.Range("D36").Value = "582.82"
$ws.Range("E36").Value = "  +0.81%  "

# Row 37 - Cosmos
$ws.Range("D37").Value = "11.02"
$ws.Range("E37").Value = "  -0.95%  "

# Row 38 - Hedera
$ws.Range("D38").Value = "0.106"
$ws.Range("E38").Value = "  +2.34%  "

# Row 39 - OKB
$ws.Range("D39").Value = "58.07"
$ws.Range("E39").Value = "  -1.67%  "

# Row 40 - FirstDigitalUSD
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.03%  "

# Row 41 - Maker
$ws.Range("D41").Value = "3.577.03"
$ws.Range("E41").Value = "  -0.72%  "

# Rows 42/43 - Kaspa and VeChain swapped ranking order
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.141"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0452"
$ws.Range("E43").Value = "  +5.59%  "

# Row 44 - TheGraph
$ws.Range("D44").Value = "0.343"
$ws.Range("E44").Value = "  +0.71%  "

# Row 45 - InjectiveProtocol
$ws.Range("D45").Value = "34.67"
$ws.Range("E45").Value = "  -2.54%  "

# Row 46 - PEPE
$ws.Range("D46").Value = "0.0₃0740"

# Row 47 - Fetch.AI
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  -1.68%  "

# Row 48 - ThetaToken
$ws.Range("D48").Value = "2.90"
$ws.Range("E48").Value = "  +6.80%  "

# Row 49 - Stellar
$ws.Range("D49").Value = "0.132"
$ws.Range("E49").Value = "  +1.57%  "

# Row 50 - Monero
$ws.Range("D50").Value = "135.54"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51 - LidoDAOToken
$ws.Range("D51").Value = "2.92"
$ws.Range("E51").Value = "  +6.37%  "
